$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was updated from
# 45202 (2023-10-03) to 45203 (2023-10-04) for every data row (rows 2-171).
for ($row = 2; $row -le 171; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
